# C1--C2-and-C3-PowerPoint.pptx edit
#
# The authoritative diff for this commit touches exactly one visible,
# user-reachable piece of presentation content: the table on slide 16
# (the "Google Shape;213;p29" graphic frame) is switched from the
# deck's custom table style ("Table_0", {ACDB1803-4905-4A67-9CC9-030483EBB864})
# to a different table style, {FA0EC0C0-C79A-42FE-92B5-728843DFD115}
# (one of PowerPoint's built-in table-style gallery entries). This is
# exactly what happens when the author reselects a table style from the
# Table Styles gallery on the Table Design ribbon tab.
#
# PowerPoint's Table object does not allow `Table.Style` to be assigned
# directly (it throws: "Table styles cannot be assigned through a
# property - call Table.ApplyStyle(...) instead"), so we use
# Table.ApplyStyle with the target style's GUID, matching what the
# Table Styles gallery does under the hood.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(16)
$shp = $s.Shapes.Item(3)
$tbl = $shp.Table
$tbl.ApplyStyle("{FA0EC0C0-C79A-42FE-92B5-728843DFD115}")
